$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that must retain their textual (string) representation even though the
# new value looks like a plain decimal number (single "." ) need an explicit
# text number-format applied before the value is written, otherwise Excel would
# silently coerce them into a numeric value (e.g. "43.00" -> 43).

# ---- Row 2 ----
$ws.Range("D2").Value = '92.753.64'
$ws.Range("E2").Value = '  -1.94%  '
# ---- Row 3 ----
$ws.Range("D3").Value = '3.410.50'
$ws.Range("E3").Value = '  -0.52%  '
# ---- Row 4 ----
$ws.Range("E4").Value = '  -0.01%  '
# ---- Row 5 ----
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '230.07'
$ws.Range("E5").Value = '  -3.35%  '
# ---- Row 6 ----
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '619.14'
$ws.Range("E6").Value = '  -3.77%  '
# ---- Row 7 ----
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.37'
$ws.Range("E7").Value = '  -5.53%  '
# ---- Row 8 ----
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.391'
$ws.Range("E8").Value = '  -3.83%  '
# ---- Row 9 ----
$ws.Range("E9").Value = '  +0.07%  '
# ---- Row 10 ----
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.961'
$ws.Range("E10").Value = '  -2.53%  '
# ---- Row 11 ----
$ws.Range("D11").Value = '3.409.46'
$ws.Range("E11").Value = '  -0.55%  '
# ---- Row 12 ----
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.00'
$ws.Range("E12").Value = '  +1.95%  '
# ---- Row 13 ----
$ws.Range("E13").Value = '  -1.58%  '
# ---- Row 14 ----
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.22'
$ws.Range("E14").Value = '  -0.61%  '
# ---- Row 15 ----
$ws.Range("D15").Value = '4.042.91'
$ws.Range("E15").Value = '  -0.64%  '
# ---- Row 16 ----
$ws.Range("D16").Value = '92.624.93'
$ws.Range("E16").Value = '  -1.80%  '
# ---- Row 17 ----
$ws.Range("E17").Value = '  -2.96%  '
# ---- Row 18 ----
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.09'
$ws.Range("E18").Value = '  -4.02%  '
# ---- Row 19 ----
$ws.Range("D19").Value = '3.407.06'
$ws.Range("E19").Value = '  -0.78%  '
# ---- Row 20 ----
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.93'
$ws.Range("E20").Value = '  +1.86%  '
# ---- Row 21 ----
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.46'
$ws.Range("E21").Value = '  -1.19%  '
# ---- Row 22 ----
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '495.94'
$ws.Range("E22").Value = '  -1.19%  '
# ---- Row 23 ----
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.31'
$ws.Range("E23").Value = '  +2.22%  '
# ---- Row 24 ----
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.434'
$ws.Range("E24").Value = '  -13.33%  '
# ---- Row 25 ----
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.52'
$ws.Range("E25").Value = '  -0.54%  '
# ---- Row 26 ----
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000184'
$ws.Range("E26").Value = '  -5.04%  '
# ---- Row 27 ----
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '94.39'
$ws.Range("E27").Value = '  +0.00%  '
# ---- Row 28 ----
$ws.Range("B28").Value = 'Aptos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.89'
$ws.Range("E28").Value = '  -0.96%  '
# ---- Row 29 ----
$ws.Range("B29").Value = 'WrappedeETH'
$ws.Range("C29").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D29").Value = '3.588.69'
$ws.Range("E29").Value = '  -0.64%  '
# ---- Row 30 ----
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '11.29'
$ws.Range("E30").Value = '  -4.47%  '
# ---- Row 31 ----
$ws.Range("B31").Value = 'Dai'
$ws.Range("C31").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").Value = '  +0.10%  '
# ---- Row 32 ----
$ws.Range("E32").Value = '  -1.77%  '
# ---- Row 33 ----
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.135'
$ws.Range("E33").Value = '  -3.40%  '
# ---- Row 34 ----
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.997'
$ws.Range("E34").Value = '  -0.26%  '
# ---- Row 35 ----
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.172'
$ws.Range("E35").Value = '  -4.70%  '
# ---- Row 36 ----
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '29.86'
$ws.Range("E36").Value = '  +0.16%  '
# ---- Row 37 ----
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.538'
$ws.Range("E37").Value = '  -2.76%  '
# ---- Row 38 ----
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '550.66'
$ws.Range("E38").Value = '  -1.52%  '
# ---- Row 39 ----
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.44'
$ws.Range("E39").Value = '  -3.42%  '
# ---- Row 40 ----
$ws.Range("E40").Value = '  -0.04%  '
# ---- Row 41 ----
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.39'
$ws.Range("E41").Value = '  -4.78%  '
# ---- Row 42 ----
$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.918'
$ws.Range("E42").Value = '  +0.68%  '
# ---- Row 43 ----
$ws.Range("B43").Value = 'Kaspa'
$ws.Range("C43").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.149'
$ws.Range("E43").Value = '  -1.54%  '
# ---- Row 44 ----
$ws.Range("B44").Value = 'WhiteBITCoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '23.65'
$ws.Range("E44").Value = '  -1.77%  '
# ---- Row 45 ----
$ws.Range("B45").Value = 'ImmutableX'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.71'
$ws.Range("E45").Value = '  -0.83%  '
# ---- Row 46 ----
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.66'
$ws.Range("E46").Value = '  -1.67%  '
# ---- Row 47 ----
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.45'
$ws.Range("E47").Value = '  -4.48%  '
# ---- Row 48 ----
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0405'
$ws.Range("E48").Value = '  -2.08%  '
# ---- Row 49 ----
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '53.22'
$ws.Range("E49").Value = '  -3.18%  '
# ---- Row 50 ----
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.10'
$ws.Range("E50").Value = '  -4.09%  '
# ---- Row 51 ----
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.97'
$ws.Range("E51").Value = '  -1.61%  '
